$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while keeping it stored as TEXT
# (matching the workbook's convention of storing all data-table cells,
# including numeric-looking ones, as strings) and without leaving any
# stray number-format style behind on the cell.
function Set-TextValue {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Value
    )
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# --- Column D (Price) updates for rows whose coin/link/rank stayed the same ---
Set-TextValue 2  4 "251.33"
Set-TextValue 3  4 "23.01"
Set-TextValue 4  4 "5.489"
Set-TextValue 5  4 "0.05661"
Set-TextValue 6  4 "3.441"
Set-TextValue 7  4 "6.406"
Set-TextValue 8  4 "0.8162"
Set-TextValue 9  4 "0.9323"
Set-TextValue 10 4 "0.1439"
Set-TextValue 11 4 "0.07489"
Set-TextValue 13 4 "0.03082"
Set-TextValue 14 4 "0.09339"
Set-TextValue 15 4 "3.552"
Set-TextValue 16 4 "0.001609"
Set-TextValue 17 4 "0.04758"

# Row 18 ("One"): price update + "Worstin24h" marker appended to the rank cell
Set-TextValue 18 4 "0.0005788"
Set-TextValue 18 5 "17OneONEWorstin24h"

Set-TextValue 19 4 "0.006365"
Set-TextValue 20 4 "0.005024"
Set-TextValue 21 4 "0.001032"
Set-TextValue 22 4 "0.0001499"
Set-TextValue 23 4 "3.725"
Set-TextValue 24 4 "2.189"
Set-TextValue 26 4 "0.1305"
Set-TextValue 28 4 "0.0002999"
Set-TextValue 40 4 "0.04012"

# --- Rows 41-43: coins rotate position (KickToken, BKEXToken, CEJI) and the
#     "Worstin24h" marker moves away from row 43 (now carried by row 18) ---
Set-TextValue 41 2 "KickToken"
Set-TextValue 41 3 "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 41 4 "0.006955"
Set-TextValue 41 5 "40KickTokenKICK"

Set-TextValue 42 2 "BKEXToken"
Set-TextValue 42 3 "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1072"
Set-TextValue 42 5 "41BKEXTokenBKK"

Set-TextValue 43 2 "CEJI"
Set-TextValue 43 3 "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 43 4 "0.002767"
Set-TextValue 43 5 "42CEJICEJI"

Set-TextValue 44 4 "0.007885"
Set-TextValue 45 4 "0.00005571"
Set-TextValue 48 4 "0.4998"
